# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2210
#   *_new -> *_FV2304
# then wrap the data range in an Excel Table ("Table1") and freeze the
# header row (row 1) in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1, columns A:U) -----------------------
$ws.Range("A1").Value2 = "Segmentname_FV2210"
$ws.Range("B1").Value2 = "Segmentgruppe_FV2210"
$ws.Range("C1").Value2 = "Segment_FV2210"
$ws.Range("D1").Value2 = "Datenelement_FV2210"
$ws.Range("E1").Value2 = "Segment ID_FV2210"
$ws.Range("F1").Value2 = "Code_FV2210"
$ws.Range("G1").Value2 = "Qualifier_FV2210"
$ws.Range("H1").Value2 = "Beschreibung_FV2210"
$ws.Range("I1").Value2 = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value2 = "Bedingung_FV2210"
$ws.Range("K1").Value2 = "diff"
$ws.Range("L1").Value2 = "Segmentname_FV2304"
$ws.Range("M1").Value2 = "Segmentgruppe_FV2304"
$ws.Range("N1").Value2 = "Segment_FV2304"
$ws.Range("O1").Value2 = "Datenelement_FV2304"
$ws.Range("P1").Value2 = "Segment ID_FV2304"
$ws.Range("Q1").Value2 = "Code_FV2304"
$ws.Range("R1").Value2 = "Qualifier_FV2304"
$ws.Range("S1").Value2 = "Beschreibung_FV2304"
$ws.Range("T1").Value2 = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value2 = "Bedingung_FV2304"

# --- 2. Turn the data range into an Excel Table (ListObject) -------------
$dataRange = $ws.Range("A1:U82")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split below row 1, pane = bottomLeft) -----
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
